# Apply edit: add averaged-intensity rows for newly run spiral sampling schemes.
# Re-orders the "Gaussian-Quadrature" row up (directly after "Ring Perpendicular to TD"),
# inserts the three new Spiral-* rows right after it, and shifts the remaining rotation /
# hex-grid scheme rows down accordingly (matches the HW10 notebook re-run output).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 17:19 need the same style as the existing id/row-index column (A10:A16).
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.636273015090748
$ws.Range("D10").Value = 2.044356414943652
$ws.Range("E10").Value = 1.021648733139336
$ws.Range("F10").Value = 1.636273015090748
$ws.Range("G10").Value = 0.7430223234523844
$ws.Range("H10").Value = 2.018785056086185
$ws.Range("I10").Value = 0.7717809265185523
$ws.Range("J10").Value = 2.044356414943652
$ws.Range("K10").Value = 1.533002574041494
$ws.Range("L10").Value = 1.584637794566121
$ws.Range("M10").Value = 1.372644411538476

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.6810375720588527
$ws.Range("D11").Value = 1.870249251876467
$ws.Range("E11").Value = 1.191963577177698
$ws.Range("F11").Value = 0.6810375720588527
$ws.Range("G11").Value = 0.661219373489575
$ws.Range("H11").Value = 2.700648528094625
$ws.Range("I11").Value = 0.8612670816559981
$ws.Range("J11").Value = 1.870249251876467
$ws.Range("K11").Value = 1.531106414527083
$ws.Range("L11").Value = 1.106071993292968
$ws.Range("M11").Value = 1.327730897392203

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.6786364579681472
$ws.Range("D12").Value = 1.837125260215813
$ws.Range("E12").Value = 1.194931837644306
$ws.Range("F12").Value = 0.6786364579681472
$ws.Range("G12").Value = 0.6575035489597678
$ws.Range("H12").Value = 2.708664893745735
$ws.Range("I12").Value = 0.8644896137906268
$ws.Range("J12").Value = 1.837125260215813
$ws.Range("K12").Value = 1.51602854893006
$ws.Range("L12").Value = 1.097332503449103
$ws.Range("M12").Value = 1.323558602054066

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.6763570335211739
$ws.Range("D13").Value = 1.865759699735092
$ws.Range("E13").Value = 1.192374955774677
$ws.Range("F13").Value = 0.6763570335211739
$ws.Range("G13").Value = 0.6601352927874836
$ws.Range("H13").Value = 2.703125683911524
$ws.Range("I13").Value = 0.8628365725463188
$ws.Range("J13").Value = 1.865759699735092
$ws.Range("K13").Value = 1.529067327754884
$ws.Range("L13").Value = 1.102712180638029
$ws.Range("M13").Value = 1.326764873046045

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 4.572532000000002
$ws.Range("D14").Value = 0.0006359999999999099
$ws.Range("E14").Value = 0.7017760000000021
$ws.Range("F14").Value = 4.572532000000002
$ws.Range("G14").Value = 0.04889599999999918
$ws.Range("H14").Value = 1.932612000000004
$ws.Range("I14").Value = 0.7890480000000004
$ws.Range("J14").Value = 0.0006359999999999099
$ws.Range("K14").Value = 0.351206000000001
$ws.Range("L14").Value = 2.461869000000002
$ws.Range("M14").Value = 1.340916666666668

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 7.651412500000028
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.04
$ws.Range("F15").Value = 7.651412500000028
$ws.Range("G15").Value = 0.06
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.4798625
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0.02
$ws.Range("L15").Value = 3.835706250000014
$ws.Range("M15").Value = 1.371879166666671

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 4.726227978342382
$ws.Range("D16").Value = 0.4302477689856002
$ws.Range("E16").Value = 0.4403918233599986
$ws.Range("F16").Value = 4.726227978342382
$ws.Range("G16").Value = 0.4606822620160008
$ws.Range("H16").Value = 0.424115874406405
$ws.Range("I16").Value = 0.7284485060607999
$ws.Range("J16").Value = 0.4302477689856002
$ws.Range("K16").Value = 0.4353197961727994
$ws.Range("L16").Value = 2.580773887257591
$ws.Range("M16").Value = 1.201685702195198

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9872582869970196
$ws.Range("D17").Value = 1.011564850288132
$ws.Range("E17").Value = 0.9878561293054662
$ws.Range("F17").Value = 0.9872582869970196
$ws.Range("G17").Value = 0.9941458784302897
$ws.Range("H17").Value = 1.014647477241821
$ws.Range("I17").Value = 0.9933815826730794
$ws.Range("J17").Value = 1.011564850288132
$ws.Range("K17").Value = 0.9997104897967991
$ws.Range("L17").Value = 0.9934843883969093
$ws.Range("M17").Value = 0.9981423674893014

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.263421914281793
$ws.Range("D18").Value = 0.822653261232401
$ws.Range("E18").Value = 1.096786058537326
$ws.Range("F18").Value = 1.263421914281793
$ws.Range("G18").Value = 0.9588729589307317
$ws.Range("H18").Value = 0.6807956322691635
$ws.Range("I18").Value = 0.9946895064433229
$ws.Range("J18").Value = 0.822653261232401
$ws.Range("K18").Value = 0.9597196598848634
$ws.Range("L18").Value = 1.111570787083328
$ws.Range("M18").Value = 0.9695365552824563

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.046630165138998
$ws.Range("D19").Value = 1.259877008313853
$ws.Range("E19").Value = 0.9267081334668097
$ws.Range("F19").Value = 1.046630165138998
$ws.Range("G19").Value = 1.067605000113277
$ws.Range("H19").Value = 0.9955707065057025
$ws.Range("I19").Value = 0.9472715231865726
$ws.Range("J19").Value = 1.259877008313853
$ws.Range("K19").Value = 1.093292570890331
$ws.Range("L19").Value = 1.069961368014665
$ws.Range("M19").Value = 1.040610422787536
